$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.4034583568573
$ws.Range("B1").Value = 1.504149675369263
$ws.Range("C1").Value = 1.687366724014282
$ws.Range("D1").Value = 2.627858877182007
$ws.Range("E1").Value = 4.607229232788086
